$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Roles y responsabilidades" column (D) had been left blank for every
# deliverable while the consistency check was being finished; fill it in
# with "si" for every reviewed row (5 through 28), matching column E/F.
$ws.Range("D5:D28").Value = "si"

# Leave the selection where the reviewer ended up after finishing the column.
$ws.Range("D29").Select() | Out-Null
